$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated 2D training schedule data (rows 2-6, columns A-J)
$data = @(
    @(1, 2, 4, 7, 5, 5, 1, 12, 5, "train_dim2_1"),
    @(2, 0, 4, 4, 6, 4, 2, 23, 5, "train_dim2_1"),
    @(3, 4, 0, 5, 5, 1, 5, 56, 5, "train_dim2_1"),
    @(4, 3, 3, 6, 6, 3, 3, 34, 5, "train_dim2_1"),
    @(5, 1, 1, 3, 5, 2, 4, 45, 5, "train_dim2_1")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $values[$j]
    }
}

# The saved workbook no longer pins the cursor on A3; reset it to the default A1 cell
$ws.Range("A1").Select()
